# Updates "想去人数" (want-to-go count, column F) figures on each sheet to
# reflect freshly generated output (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value  = 1172
$ws.Range("F8").Value  = 1076
$ws.Range("F10").Value = 375
$ws.Range("F15").Value = 44
$ws.Range("F18").Value = 556
$ws.Range("F20").Value = 5719
$ws.Range("F22").Value = 1609
$ws.Range("F25").Value = 32
$ws.Range("F26").Value = 5304
$ws.Range("F27").Value = 5304
$ws.Range("F30").Value = 1541

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 33
$ws.Range("F5").Value = 159

# Sheet 3: 本地生活 (Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 9414

# Sheet 4: 全部类型 (All types - combined view)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value  = 9414
$ws.Range("F7").Value  = 1172
$ws.Range("F11").Value = 1076
$ws.Range("F12").Value = 375
$ws.Range("F16").Value = 44
$ws.Range("F22").Value = 5719
$ws.Range("F24").Value = 1609
$ws.Range("F30").Value = 5304
$ws.Range("F31").Value = 5304
$ws.Range("F34").Value = 1541
